# Daily attendance processing - 2026-01-14 19:59:14
# Reorders the "Recorded By" (column G) values for specific exact matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row so we cover the whole data range.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

# Map of exact old values -> new values (reordering the comma-separated entries).
$replacements = @{
    "dnasr281@gmail.com, System"               = "System, dnasr281@gmail.com"
    "System, backup@backdoor.com"               = "backup@backdoor.com, System"
    "system, System, backup@backdoor.com"       = "backup@backdoor.com, System, system"
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
